{"js": "// Remove the trailing space after the period that ends each top-level\n// heading's run (e.g. \"2.\\tConfidential Information. \" -> \"...Information.\").\n// Only the final run of each matching heading paragraph is touched; every\n// other run (and its formatting) is left untouched.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const paragraph = paragraphs.items[i];\n  // Heading paragraphs in this document end with a lone \". \" run right\n  // after the bold heading text (e.g. \"2.\\tConfidential Information. \").\n  if (paragraph.text.endsWith(\". \")) {\n    targets.push(paragraph);\n  }\n}\n\nfor (const paragraph of targets) {\n  // Scope the search to this paragraph only, and grab the last match so we\n  // edit the trailing \". \" run at the end of the paragraph (not some other\n  // \". \" substring earlier in the text).\n  const found = paragraph.getRange().search(\". \", { matchCase: true });\n  found.load(\"items/text\");\n  await context.sync();\n\n  if (found.items.length === 0) {\n    continue;\n  }\n\n  const lastMatch = found.items[found.items.length - 1];\n  // Replacing just this sub-range keeps it inside the existing trailing\n  // run, so the edit becomes a plain text-content change on that run\n  // rather than a restructuring of the paragraph's runs.\n  lastMatch.insertText(\".\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing space after the period that ends each top-level\n# heading's run (e.g. \"2.`tConfidential Information. \" -> \"...Information.\").\n# Only the final character (the trailing space) of each matching heading\n# paragraph is deleted, so the rest of the paragraph's runs/formatting stay\n# untouched.\n\n$d = $word.ActiveDocument\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    $text = $r.Text\n\n    # Heading paragraphs in this document end with a lone \". \" run right\n    # after the bold heading text (e.g. \"2.`tConfidential Information. \").\n    # $text includes the trailing paragraph mark, so check for \". `r\".\n    if ($text -like \"*. `r\") {\n        $count = $r.Characters.Count\n        # Character $count is the paragraph mark; $count - 1 is the space\n        # right before it that needs to be removed.\n        $spaceChar = $r.Characters.Item($count - 1)\n        $spaceChar.Delete()\n    }\n}\n"}
